# Apply cryptos list price/volume/coin updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Cell='D2'; Value='60.854.24'}
    @{Cell='E2'; Value='  +2.70%  '}
    @{Cell='D3'; Value='2.631.29'}
    @{Cell='E3'; Value='  +2.05%  '}
    @{Cell='E4'; Value='  +0.10%  '}
    @{Cell='D5'; Value='''570.23'}
    @{Cell='E5'; Value='  -0.44%  '}
    @{Cell='D6'; Value='''142.97'}
    @{Cell='E6'; Value='  -0.09%  '}
    @{Cell='D7'; Value='''0.998'}
    @{Cell='E7'; Value='  -0.11%  '}
    @{Cell='E8'; Value='  +0.90%  '}
    @{Cell='D9'; Value='2.630.36'}
    @{Cell='E9'; Value='  +1.80%  '}
    @{Cell='E10'; Value='  -3.08%  '}
    @{Cell='E11'; Value='  +2.64%  '}
    @{Cell='D12'; Value='''0.153'}
    @{Cell='E12'; Value='  -3.66%  '}
    @{Cell='E13'; Value='  +6.94%  '}
    @{Cell='D14'; Value='3.071.36'}
    @{Cell='E14'; Value='  +1.34%  '}
    @{Cell='D15'; Value='60.855.52'}
    @{Cell='E15'; Value='  +2.69%  '}
    @{Cell='D16'; Value='''23.54'}
    @{Cell='E16'; Value='  +5.01%  '}
    @{Cell='E17'; Value='  +3.05%  '}
    @{Cell='D18'; Value='2.619.94'}
    @{Cell='E18'; Value='  +1.45%  '}
    @{Cell='D19'; Value='''11.22'}
    @{Cell='E19'; Value='  +9.29%  '}
    @{Cell='D20'; Value='''4.66'}
    @{Cell='E20'; Value='  +2.82%  '}
    @{Cell='D21'; Value='''349.07'}
    @{Cell='E21'; Value='  +3.15%  '}
    @{Cell='E22'; Value='  +12.83%  '}
    @{Cell='E23'; Value='  +0.16%  '}
    @{Cell='E24'; Value='  +12.71%  '}
    @{Cell='D25'; Value='''64.49'}
    @{Cell='E25'; Value='  -0.21%  '}
    @{Cell='E26'; Value='  -0.32%  '}
    @{Cell='E27'; Value='  -0.08%  '}
    @{Cell='D28'; Value='''7.72'}
    @{Cell='E28'; Value='  +5.99%  '}
    @{Cell='D29'; Value='0.0₃0796'}
    @{Cell='E29'; Value='  +1.69%  '}
    @{Cell='D30'; Value='''1.83'}
    @{Cell='E30'; Value='  +8.44%  '}
    @{Cell='D31'; Value='''0.998'}
    @{Cell='E31'; Value='  -0.05%  '}
    @{Cell='E32'; Value='  +4.57%  '}
    @{Cell='D33'; Value='''160.59'}
    @{Cell='E33'; Value='  +1.19%  '}
    @{Cell='D34'; Value='''19.54'}
    @{Cell='E34'; Value='  +2.69%  '}
    @{Cell='E35'; Value='  +5.79%  '}
    @{Cell='D36'; Value='''0.963'}
    @{Cell='E36'; Value='  +9.72%  '}
    @{Cell='E37'; Value='  +3.64%  '}
    @{Cell='E38'; Value='  +7.12%  '}
    @{Cell='D39'; Value='''37.81'}
    @{Cell='E39'; Value='  +1.37%  '}
    @{Cell='D40'; Value='''0.851'}
    @{Cell='E40'; Value='  -2.69%  '}
    @{Cell='D41'; Value='''3.80'}
    @{Cell='E41'; Value='  +3.56%  '}
    @{Cell='D42'; Value='''298.28'}
    @{Cell='E42'; Value='  +0.70%  '}
    @{Cell='D43'; Value='''140.35'}
    @{Cell='E43'; Value='  +7.39%  '}
    @{Cell='B44'; Value='FirstDigitalUSD'}
    @{Cell='C44'; Value='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'}
    @{Cell='D44'; Value='''0.998'}
    @{Cell='E44'; Value='  -0.09%  '}
    @{Cell='B45'; Value='Stellar'}
    @{Cell='C45'; Value='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'}
    @{Cell='D45'; Value='''0.0987'}
    @{Cell='E45'; Value='  +0.89%  '}
    @{Cell='E46'; Value='  +2.11%  '}
    @{Cell='E47'; Value='  +2.16%  '}
    @{Cell='B48'; Value='EnergySwap'}
    @{Cell='C48'; Value='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'}
    @{Cell='D48'; Value='''19.60'}
    @{Cell='E48'; Value='  +2.14%  '}
    @{Cell='B49'; Value='VeChain'}
    @{Cell='C49'; Value='https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'}
    @{Cell='D49'; Value='''0.0242'}
    @{Cell='E49'; Value='  +3.36%  '}
    @{Cell='B50'; Value='WhiteBITCoin'}
    @{Cell='C50'; Value='https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'}
    @{Cell='D50'; Value='''10.71'}
    @{Cell='E50'; Value='  +0.49%  '}
    @{Cell='B51'; Value='InjectiveProtocol'}
    @{Cell='C51'; Value='https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'}
    @{Cell='D51'; Value='''19.74'}
    @{Cell='E51'; Value='  +6.29%  '}
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"